$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row was inserted above the existing row 138.
# This shifts every following row (old 138..221) down by one (to 139..222)
# and keeps the sheet's dimension in sync (A1:R221 -> A1:R222).
$ws.Rows("138:138").Insert()

# Populate the freshly inserted row with the new record's data.
$ws.Cells.Item(138, 1).Value  = 7
$ws.Cells.Item(138, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(138, 3).Value  = "Ñuble"
$ws.Cells.Item(138, 4).Value  = 44777
$ws.Cells.Item(138, 5).Value  = 16
$ws.Cells.Item(138, 6).Value  = 100112017
$ws.Cells.Item(138, 7).Value  = "Apio"
$ws.Cells.Item(138, 8).Value  = "Americana (o)"
$ws.Cells.Item(138, 9).Value  = "Primera"
$ws.Cells.Item(138, 10).Value = 120
$ws.Cells.Item(138, 11).Value = 9000
$ws.Cells.Item(138, 12).Value = 10000
$ws.Cells.Item(138, 13).Value = 9500
$ws.Cells.Item(138, 14).Value = "$/docena de matas"
$ws.Cells.Item(138, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(138, 16).Value = 1583
$ws.Cells.Item(138, 17).Value = 6
$ws.Cells.Item(138, 18).Value = "Hortaliza"

# Ensure the date cell keeps the workbook's date number format (style index 2),
# matching every other cell in column D.
$ws.Cells.Item(138, 4).NumberFormat = $ws.Cells.Item(139, 4).NumberFormat
